# Auto-generated: refresh the crypto price/volume columns (D, E) to match
# the latest GitHub Actions data pull, and swap the Aptos / WEMIXTOKEN rows
# (37 <-> 38) which changed rank order in the source feed.
#
# Cells in column D sometimes hold values that are valid numbers (e.g.
# "1.001"); Excel's normal text-entry parsing would silently convert those
# to floating point and drop the original text formatting (trailing zeros,
# fixed decimal places) that the source site renders. Force those specific
# cells to Text format first so the literal string is preserved exactly,
# matching how the source data is actually stored (inline text, not a
# number). Cells that are already unambiguous text (URLs, names, values with
# two '.' separators, the padded "  +/-X.XX%  " volume strings) do not need
# this and are assigned directly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.935.48'
$ws.Range("E2").Value = '  -3.21%  '
$ws.Range("D3").Value = '1.736.94'
$ws.Range("E3").Value = '  -1.29%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.001'
$ws.Range("E4").Value = '  -0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '310.51'
$ws.Range("E5").Value = '  -5.65%  '
$ws.Range("E6").Value = '  -0.02%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4972'
$ws.Range("E7").Value = '  +3.04%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3534'
$ws.Range("E8").Value = '  -0.28%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '42.46'
$ws.Range("E9").Value = '  -2.26%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07234'
$ws.Range("E10").Value = '  -4.70%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.056'
$ws.Range("E11").Value = '  -2.26%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.001'
$ws.Range("E12").Value = '  -0.03%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '19.91'
$ws.Range("E13").Value = '  -2.99%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.930'
$ws.Range("E14").Value = '  -2.43%  '
$ws.Range("D15").Value = '1.734.70'
$ws.Range("E15").Value = '  -2.02%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.837'
$ws.Range("E16").Value = '  -4.45%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '86.13'
$ws.Range("E17").Value = '  -6.75%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.00001033'
$ws.Range("E18").Value = '  -5.84%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06379'
$ws.Range("E19").Value = '  -0.74%  '
$ws.Range("E20").Value = '  -0.06%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '16.56'
$ws.Range("E21").Value = '  -1.27%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.727'
$ws.Range("E22").Value = '  -1.06%  '
$ws.Range("D23").Value = '27.001.03'
$ws.Range("E23").Value = '  -3.10%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.14'
$ws.Range("E24").Value = '  +0.33%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.047'
$ws.Range("E25").Value = '  -5.11%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '153.44'
$ws.Range("E26").Value = '  -6.19%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '19.84'
$ws.Range("E27").Value = '  -1.09%  '
$ws.Range("D28").Value = '1.934.64'
$ws.Range("E28").Value = '  -1.88%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.117'
$ws.Range("E29").Value = '  -3.93%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '120.35'
$ws.Range("E30").Value = '  -2.08%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.052'
$ws.Range("E31").Value = '  -0.40%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.09435'
$ws.Range("E32").Value = '  -0.40%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.574'
$ws.Range("E33").Value = '  -1.97%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.357'
$ws.Range("E34").Value = '  -3.52%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.05909'
$ws.Range("E35").Value = '  -1.38%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.02187'
$ws.Range("E36").Value = '  -3.23%  '
$ws.Range("B37").Value = 'Aptos'
$ws.Range("C37").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '11.00'
$ws.Range("E37").Value = '  -5.25%  '
$ws.Range("B38").Value = 'WEMIXTOKEN'
$ws.Range("C38").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.426'
$ws.Range("E38").Value = '  -0.36%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.1994'
$ws.Range("E39").Value = '  -3.30%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '4.741'
$ws.Range("E40").Value = '  -3.20%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.000'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.5989'
$ws.Range("E42").Value = '  -2.59%  '
$ws.Range("E43").Value = '  -6.27%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '7.428'
$ws.Range("E44").Value = '  -3.54%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '12.85'
$ws.Range("E45").Value = '  -2.02%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.574'
$ws.Range("E46").Value = '  -4.38%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5621'
$ws.Range("E47").Value = '  -2.86%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '119.78'
$ws.Range("E48").Value = '  -3.19%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.849'
$ws.Range("E49").Value = '  -3.74%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06658'
$ws.Range("E50").Value = '  -1.91%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.096'
$ws.Range("E51").Value = '  -4.87%  '
